$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long pandas-repr text blocks with the placeholder text.
$rowsToReplace = @(6, 14, 23, 32, 53, 60, 62, 74, 79, 88, 89, 91, 97, 100, 110)
foreach ($r in $rowsToReplace) {
    $ws.Cells.Item($r, 1).Value = "no need to revise"
}

# Drop the trailing 4 rows (117-120) entirely, shrinking the used range to A1:A116.
$ws.Range("A117:A120").EntireRow.Delete()
